# Reorder "Recorded By" names in column G so that "System" (capitalized)
# is listed first among the comma-separated names, matching the upstream
# commit that synced attendance_reports data.
#
# Two observed patterns in this sheet:
#   "dnasr281@gmail.com, System"                 -> "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com"        -> "System, system, backup@backdoor.com"
#
# In both cases the token "System" (exact case) is moved to the front of
# the comma-separated list, preserving the relative order of the other
# tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value -notmatch ",") { continue }

    $parts = $value -split ", "
    if ($parts -notcontains "System") { continue }

    # NOTE: string comparison operators (-eq/-ne/-ceq/-cne) in this host
    # behave case-insensitively even with the "c" prefix, so avoid relying
    # on them to distinguish "System" vs "system". [array]::IndexOf() on a
    # string array does a case-sensitive match, which is what we need here.
    $idx = [array]::IndexOf($parts, "System")
    if ($idx -eq 0) { continue }  # already first, nothing to do

    # Build the remaining tokens (all but the "System" one) via an explicit
    # loop rather than range slicing, since $parts[($idx+1)..($parts.Length-1)]
    # misbehaves when the computed range would be empty (idx is the last
    # element) in this host's PowerShell-subset.
    $rest = @()
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($i -ne $idx) { $rest += $parts[$i] }
    }

    $reordered = @("System") + $rest
    $newValue = $reordered -join ", "

    $cell.Value2 = $newValue
}
